$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2:D2").Copy()
$ws.Range("A3:D8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Pass 1: project/task/date columns (keeps shared-string order: dates before remarks)
$ws.Cells.Item(2, 1).Value = "BBBY-PMall"
$ws.Cells.Item(2, 2).Value = "Test Execution/Defect Documentation/Defect Retest/Testing"
$ws.Cells.Item(2, 3).Value = "10/21/2019"
$ws.Cells.Item(3, 1).Value = "BBBY-PMall"
$ws.Cells.Item(3, 2).Value = "Test Execution/Defect Documentation/Defect Retest/Testing"
$ws.Cells.Item(3, 3).Value = "10/22/2019"
$ws.Cells.Item(4, 1).Value = "BBBY-PMall"
$ws.Cells.Item(4, 2).Value = "Test Execution/Defect Documentation/Defect Retest/Testing"
$ws.Cells.Item(4, 3).Value = "10/23/2019"
$ws.Cells.Item(5, 1).Value = "BBBY-PMall"
$ws.Cells.Item(5, 2).Value = "Test Execution/Defect Documentation/Defect Retest/Testing"
$ws.Cells.Item(5, 3).Value = "10/24/2019"
$ws.Cells.Item(6, 1).Value = "BBBY-PMall"
$ws.Cells.Item(6, 2).Value = "Test Execution/Defect Documentation/Defect Retest/Testing"
$ws.Cells.Item(6, 3).Value = "10/25/2019"
$ws.Cells.Item(7, 1).Value = "BBBY-PMall"
$ws.Cells.Item(7, 2).Value = "Test Execution/Defect Documentation/Defect Retest/Testing"
$ws.Cells.Item(7, 3).Value = "10/29/2019"
$ws.Cells.Item(8, 1).Value = "BBBY-PMall"
$ws.Cells.Item(8, 2).Value = "Test Execution/Defect Documentation/Defect Retest/Testing"
$ws.Cells.Item(8, 3).Value = "10/30/2019"

# Pass 2: remarks column, filled from row 8 (newest date) up to row 2 (oldest date)
# to reproduce the original shared-string insertion order.
$remark8 = @'
1.	Worked on creation of test cases for October release. Please find attached sheet for October release with test cases and defects created so far.
2.	Reviewed the automation test code for understanding.
3.	Please review all defects added to October release.
4.	Need to discuss the Varvy SEO Tool tasks, before moving on further, regarding what issues are to be reported in what format and also about the scope of verification.
'@
$ws.Cells.Item(8, 4).Value = $remark8
$remark7 = @'
1.	Worked on creation of test cases for October release. Please find attached sheet for October release with test cases and defects created so far.
2.	Created defect #8572 and added the same to October release defects.
3.	Please review all defects added to October release.
4.	Need to discuss the Varvy SEO Tool tasks, before moving on further, regarding what issues are to be reported in what format and also about the scope of verification.
'@
$ws.Cells.Item(7, 4).Value = $remark7
$remark6 = @'
1.	Worked on creation of test cases for October release tickets discussed.
2.	Reviewed the automation test automation cases added by Noah in smartsheet.
3.	Please review all defects added to October release.
4.	Created defect for blank shipping address issue in PMall admin assign d the same to you.
5.	Need to discuss the Varvy SEO Tool tasks, before moving on further, regarding what issues are to be reported in what format and also about the scope of verification.
'@
$ws.Cells.Item(6, 4).Value = $remark6
$remark5 = @'
1.	Reviewed all defects added to October release. Prioritized and reordered all the defects after review as per my understanding. Also, closed few of them which were working fine but were not assigned to me. However, I have closed them and marked them as Lived in the October release defects in smartsheet.
2.	Working on  understanding the automation scripts code for desktop and mobile site, further. Not able to pull the code in the existing developer branch as some SSL issue is showing up which has been confirmed by Noah too during the automation meeting on call. 
3.	Reviewed the mockups and other requirements discussed yesterday which are part of October release.
4.	Attended daily automation meeting with Noah. 
5.	Need to discuss the Varvy SEO Tool tasks, before moving on further, regarding what issues are to be reported in what format and also about the scope of verification.
'@
$ws.Cells.Item(5, 4).Value = $remark5
$remark4 = @'
1.	All changes have been incorporated in smart sheet “eCommerce automation” for both Desktop and Mobile site scenarios. Kindly review and suggest.
2.	Reviewed the Power shell scripts on AutoQA2 server, for scheduling the automation tasks. Also, reviewing the automation scripts code for desktop and mobile site in detail. Not able to pull the code in the existing developer branch as some SSL issue is showing up which has been confirmed by Noah too during the automation meeting on call. 
3.	Attended daily automation meeting with Noah. 
4.	Attended daily status meeting to discuss the requirements of October release.
5.	Need to discuss the Varvy SEO Tool tasks, before moving on further, regarding what issues are to be reported in what format and also about the scope of verification.
'@
$ws.Cells.Item(4, 4).Value = $remark4
$remark3 = @'
1.	Updated Test scenarios for automation in smartsheet “eCommerce automation” for both Desktop and Mobile site scenarios in smartsheet. Please review and suggest. Also, please review the backlog section in the  smartsheet for further task assignment in Automation.
2.	Reviewed the Power shell scripts on AutoQA2 server, for scheduling the automation tasks. Also, reviewing the automation scripts code for desktop and mobile site in detail. 
3.	Need to discuss the Varvy SEO Tool tasks, before moving on further, regarding what issues are to be reported in what format and also about the scope of verification.
4.	Reviewed the requirements of the October release tickets.
'@
$ws.Cells.Item(3, 4).Value = $remark3
$remark2 = @'
1.	Updated all mobile site scenarios and collaborated them with Desktop ones in the required format in automation test scenarios smartsheet. Added few comments for Noah to update. Please review and suggest.
2.	Verified few PMall pages using the Varvy  SEO inspector tool as per the assigned ticket #8394. Need to discuss it before moving on what issues to be reported in what format and also about the scope of verification.
3.	Attended daily automation discussion meeting with Noah. Discussed the changes made in Automation test cases in smartsheet.
4.	Reviewed the requirements of the October release tickets.
'@
$ws.Cells.Item(2, 4).Value = $remark2

# Row heights
$ws.Rows.Item(2).RowHeight = 135
$ws.Rows.Item(3).RowHeight = 135
$ws.Rows.Item(4).RowHeight = 150
$ws.Rows.Item(5).RowHeight = 195
$ws.Rows.Item(6).RowHeight = 105
$ws.Rows.Item(7).RowHeight = 90
$ws.Rows.Item(8).RowHeight = 90

# Sheet view: scroll position + selection
$ws.Application.Goto($ws.Range("A6"), $true) | Out-Null
$ws.Range("D6").Select() | Out-Null